# Add 2022-Q3 data
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 right after the header,
#    pushing the existing quarters down by one row.
# 2) Add a new worksheet "2022-Q3" (positioned right after "总计", before "2022-Q1")
#    holding the per-fund holdings detail for that quarter.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: update the "总计" overview sheet
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Column A is just a 0-based row index (0,1,2,3,4 for rows 2..6) and is left
# untouched except for the brand-new row 7 (index 5). Only the B:D (quarter
# label / count / value) columns actually shift down by one row.
# Shift existing data rows (2..6) down to (3..7), working bottom-up so we
# never clobber a row before it has been copied.
$summary.Range("B6:D6").Copy($summary.Range("B7:D7"))
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))

# Extend column A's index sequence to the new row 7.
$summary.Range("A6").Copy($summary.Range("A7"))
$summary.Range("A7").Value = 5

# Write the new 2022-Q3 row into the now-vacated row 2 (A2 already holds 0).
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.15

# ------------------------------------------------------------------
# Step 2: create the "2022-Q3" detail sheet
# ------------------------------------------------------------------
# Duplicate the existing "2022-Q1" sheet (same column layout/formatting)
# and drop it in right before "2022-Q1" -> ends up right after "总计".
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($template)
$new = $wb.Worksheets.Item("2022-Q1 (2)")
$new.Name = "2022-Q3"

# The duplicated sheet has 2 data rows (rows 2-3); we need 3 (rows 2-4), so
# clone the formatting of row 3 down into row 4 before overwriting values.
$new.Range("A3:H3").Copy($new.Range("A4:H4"))

# Row 2
$new.Range("A2").Value = 0
$new.Range("B2").Value = "'001672"
$new.Range("C2").Value = "国寿安保智慧生活股票"
$new.Range("D2").Value = "'4.26"
$new.Range("E2").Value = "'86.24"
$new.Range("F2").Value = "'2.82"
$new.Range("G2").Value = "'0.1201"
$new.Range("H2").Value = 10

# Row 3
$new.Range("A3").Value = 1
$new.Range("B3").Value = "'010821"
$new.Range("C3").Value = "东方红多元策略混合B"
$new.Range("D3").Value = "'0.59"
$new.Range("E3").Value = "'92.90"
$new.Range("F3").Value = "'3.22"
$new.Range("G3").Value = "'0.0190"
$new.Range("H3").Value = 8

# Row 4
$new.Range("A4").Value = 2
$new.Range("B4").Value = "'910017"
$new.Range("C4").Value = "东方红多元策略混合A"
$new.Range("D4").Value = "'0.41"
$new.Range("E4").Value = "'92.90"
$new.Range("F4").Value = "'3.22"
$new.Range("G4").Value = "'0.0132"
$new.Range("H4").Value = 8
